$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.583.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.876.56"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.55%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.76"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5102"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3924"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08392"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.112"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.73"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.271"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.880.08"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.75%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.277"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001106"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.42"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06723"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.74"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.68%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.968"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.619.03"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.14"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.88%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.245"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.104.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.44"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.82"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.98%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.372"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.80"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1054"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.18%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.807"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.620"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02458"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06544"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2189"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.33%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.266"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.196"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6475"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.76%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.070"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.57%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.76%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6071"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.98"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.695"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.035"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.219"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.45"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.09%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -6.31%  "
